$wb = $excel.ActiveWorkbook

# ALC row 7
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 10000
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

# ALC row 14
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H14").Value = 10000
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

# ALC row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1145.92
$ws.Range("I15").Value = 1145.92
$ws.Range("K15").Value = 3437.76
$ws.Range("M15").Value = -3268.76

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1078.75
$ws.Range("J112").Value = 1078.75
$ws.Range("L112").Value = 3236.25
$ws.Range("N112").Value = -5452.25

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3796.6667
$ws.Range("I116").Value = 2633.3333
$ws.Range("J116").Value = 4960
$ws.Range("K116").Value = 2633.3333
$ws.Range("L116").Value = 4960
$ws.Range("M116").Value = 808.6667000000002
$ws.Range("N116").Value = -11844

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 855.23883
$ws.Range("J129").Value = 903.7931
$ws.Range("L129").Value = 2711.3793
$ws.Range("N129").Value = -12711.3793

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4587.42
$ws.Range("I138").Value = 2291.2856
$ws.Range("J138").Value = 5480.3613
$ws.Range("K138").Value = 6873.8568
$ws.Range("L138").Value = 16441.0839
$ws.Range("M138").Value = -1733.8568
$ws.Range("N138").Value = -26721.0839

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 34451.5
$ws.Range("I32").Value = 5845.2095
$ws.Range("J32").Value = 171126
$ws.Range("K32").Value = 5845.2095
$ws.Range("L32").Value = 171126
$ws.Range("M32").Value = -5558.2095
$ws.Range("N32").Value = -171700

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2283.6667
$ws.Range("I61").Value = 1989.3334
$ws.Range("J61").Value = 3166.6667
$ws.Range("K61").Value = 1989.3334
$ws.Range("L61").Value = 3166.6667
$ws.Range("M61").Value = -1777.3334
$ws.Range("N61").Value = -3590.6667

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1879.6428
$ws.Range("I132").Value = 1567.4166
$ws.Range("J132").Value = 3753
$ws.Range("K132").Value = 4702.2498
$ws.Range("L132").Value = 11259
$ws.Range("M132").Value = -2172.2498
$ws.Range("N132").Value = -16319

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2283.6667
$ws.Range("I136").Value = 1989.3334
$ws.Range("J136").Value = 3166.6667
$ws.Range("K136").Value = 5968.0002
$ws.Range("L136").Value = 9500.000100000001
$ws.Range("M136").Value = -3418.0002
$ws.Range("N136").Value = -14600.0001

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2386.7334
$ws.Range("I134").Value = 2510.5
$ws.Range("J134").Value = 1582.25
$ws.Range("K134").Value = 7531.5
$ws.Range("L134").Value = 4746.75
$ws.Range("M134").Value = -4996.5
$ws.Range("N134").Value = -9816.75

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21206.96
$ws.Range("I31").Value = 25084.404
$ws.Range("J31").Value = 4921.7
$ws.Range("K31").Value = 25084.404
$ws.Range("L31").Value = 4921.7
$ws.Range("M31").Value = -24789.404
$ws.Range("N31").Value = -5511.7

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 21206.96
$ws.Range("I34").Value = 25084.404
$ws.Range("J34").Value = 4921.7
$ws.Range("K34").Value = 25084.404
$ws.Range("L34").Value = 4921.7
$ws.Range("M34").Value = -24882.404
$ws.Range("N34").Value = -5325.7

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3270396.2
$ws.Range("I62").Value = 11113051
$ws.Range("J62").Value = 2623.4167
$ws.Range("K62").Value = 11113051
$ws.Range("L62").Value = 2623.4167
$ws.Range("M62").Value = -11112427
$ws.Range("N62").Value = -3871.4167

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 3270396.2
$ws.Range("I65").Value = 11113051
$ws.Range("J65").Value = 2623.4167
$ws.Range("K65").Value = 55565255
$ws.Range("L65").Value = 13117.0835
$ws.Range("M65").Value = -55562135
$ws.Range("N65").Value = -19357.0835

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1958.2222
$ws.Range("I134").Value = 1958.2222
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5874.6666
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -3339.6666
$ws.Range("N134").ClearContents()

# CUL row 14
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1151
$ws.Range("I14").Value = 1151
$ws.Range("K14").Value = 3453
$ws.Range("M14").Value = -3280

# CUL row 32
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 4965
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 4965
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 14895
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -15461

# CUL row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2279.9092
$ws.Range("I34").Value = 90
$ws.Range("J34").Value = 2766.5557
$ws.Range("K34").Value = 270
$ws.Range("L34").Value = 8299.667099999999
$ws.Range("M34").Value = -186
$ws.Range("N34").Value = -8467.667099999999

# CUL row 86
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1368
$ws.Range("J86").Value = 1552.5
$ws.Range("L86").Value = 4657.5
$ws.Range("N86").Value = -7029.5

# CUL row 89
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 1368
$ws.Range("J89").Value = 1552.5
$ws.Range("L89").Value = 13972.5
$ws.Range("N89").Value = -25828.5

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 687.63635
$ws.Range("I113").Value = 667.5
$ws.Range("J113").Value = 699.1429000000001
$ws.Range("K113").Value = 2002.5
$ws.Range("L113").Value = 2097.4287
$ws.Range("M113").Value = 167.5
$ws.Range("N113").Value = -6437.4287

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 846.3
$ws.Range("I131").Value = 665.8
$ws.Range("J131").Value = 855.8
$ws.Range("K131").Value = 1997.4
$ws.Range("L131").Value = 2567.4
$ws.Range("M131").Value = 3042.6
$ws.Range("N131").Value = -12647.4

# CUL row 136
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 2400.1667
$ws.Range("I136").Value = 1333.6666
$ws.Range("J136").Value = 3466.6667
$ws.Range("K136").Value = 4000.9998
$ws.Range("L136").Value = 10400.0001
$ws.Range("M136").Value = 1099.0002
$ws.Range("N136").Value = -20600.0001

# GSM row 23
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 6000
$ws.Range("J23").Value = 6000
$ws.Range("L23").Value = 6000
$ws.Range("N23").Value = -6446

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2462.5
$ws.Range("I126").Value = 2405.25
$ws.Range("J126").Value = 2519.75
$ws.Range("K126").Value = 7215.75
$ws.Range("L126").Value = 7559.25
$ws.Range("M126").Value = -4745.75
$ws.Range("N126").Value = -12499.25

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2580.25
$ws.Range("I132").Value = 2146.4
$ws.Range("K132").Value = 6439.200000000001
$ws.Range("M132").Value = -3909.200000000001

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1687785
$ws.Range("I46").Value = 700
$ws.Range("J46").Value = 2025202
$ws.Range("K46").Value = 700
$ws.Range("L46").Value = 2025202
$ws.Range("M46").Value = -512
$ws.Range("N46").Value = -2025578

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1251.9166
$ws.Range("J93").Value = 929.3333
$ws.Range("L93").Value = 929.3333
$ws.Range("N93").Value = -3425.3333

# WVR row 41
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 26067.4
$ws.Range("J41").Value = 26067.4
$ws.Range("L41").Value = 26067.4
$ws.Range("N41").Value = -26847.4

# WVR row 45
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 9217
$ws.Range("J45").Value = 9217
$ws.Range("L45").Value = 9217
$ws.Range("N45").Value = -10199

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1984.579
$ws.Range("I122").Value = 1200.1333
$ws.Range("J122").Value = 4926.25
$ws.Range("K122").Value = 3600.3999
$ws.Range("L122").Value = 14778.75
$ws.Range("M122").Value = -1150.3999
$ws.Range("N122").Value = -19678.75
